# Rewrite the bullets under "KEY ACHIEVEMENTS AND IMPACT" as short,
# impact-focused accomplishment statements (formatted like Core
# Competencies), dropping the section from 6 bullets down to 4.
#
# Several of the old/new bullet strings share substrings with bullets
# elsewhere in the document (e.g. the "PROFESSIONAL EXPERIENCE" section
# also mentions the trigonometric-algorithm / race-coding items), so all
# Find/Replace and delete operations below are scoped to a freshly
# recomputed Range that covers only the "KEY ACHIEVEMENTS AND IMPACT" ...
# "TECHNICAL SKILLS" span, rather than operating on $d.Content. The
# section is re-located by paragraph text before each edit because Range
# start/end offsets captured earlier do not automatically track edits
# made elsewhere in the document.

$d = $word.ActiveDocument
$bullet = [char]0x2022

function Get-SectionRange {
    $startParaIdx = $null
    $endParaIdx = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text.Trim()
        if ($t -eq "KEY ACHIEVEMENTS AND IMPACT" -and $startParaIdx -eq $null) {
            $startParaIdx = $i
        }
        if ($t -eq "TECHNICAL SKILLS" -and $startParaIdx -ne $null -and $endParaIdx -eq $null) {
            $endParaIdx = $i
        }
    }
    $sectionStart = $d.Paragraphs.Item($startParaIdx).Range.Start
    $sectionEnd = $d.Paragraphs.Item($endParaIdx).Range.Start
    return $d.Range($sectionStart, $sectionEnd)
}

function Replace-InSection {
    param($oldText, $newText)
    $rng = Get-SectionRange
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

function Delete-ParaInSection {
    param($text)
    $rng = Get-SectionRange
    $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    if ($rng.Find.Found) {
        # Extend the matched range by one character so the deletion also
        # removes the trailing paragraph mark (otherwise Word leaves behind
        # an empty paragraph where the bullet used to be).
        $full = $d.Range($rng.Start, $rng.End + 1)
        $full.Delete()
    }
}

# 1) Race/demographic-classification-accuracy bullet -> trigonometric
#    boundary-estimation / mapping-cost bullet.
$old1 = $bullet + " Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%"
$new1 = $bullet + " Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
Replace-InSection $old1 $new1

# 2) Redistricting-platform-usage bullet -> $4.7M savings bullet.
$old2 = $bullet + " Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations"
$new2 = $bullet + " `$4.7M savings enabled nonprofit access"
Replace-InSection $old2 $new2

# 3) Trigonometric-algorithm / mapping-cost bullet -> demographic
#    miscoding discovery bullet.
$old3 = $bullet + " Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis"
$new3 = $bullet + " Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"
Replace-InSection $old3 $new3

# 4) Longitudinal-data-analysis-methods bullet -> accuracy-improvement bullet.
$old4 = $bullet + " Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality"
$new4 = $bullet + " 178% accuracy improvement in racial classification algorithms"
Replace-InSection $old4 $new4

# 5) & 6) Drop the two trailing bullets (ETL pipelines / cloud warehouse)
#    entirely, since the rewritten section only has four bullets.
$old5 = $bullet + " Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets"
Delete-ParaInSection $old5

$old6 = $bullet + " Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy"
Delete-ParaInSection $old6
